$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the text-shared-string "xx xxx" formatted numbers with real numeric values
$ws.Range("B2").Value = 21997
$ws.Range("C2").Value = 22200

$ws.Range("B3").Value = 4750
$ws.Range("C3").Value = 4800

$ws.Range("B4").Value = 3643
$ws.Range("C4").Value = 3685

$ws.Range("B5").Value = 3180
$ws.Range("C5").Value = 3300

$ws.Range("B6").Value = 2500
$ws.Range("C6").Value = 2600

$ws.Range("B7").Value = 2400
$ws.Range("C7").Value = 2450

$ws.Range("B8").Value = 1793
$ws.Range("C8").Value = 1791

$ws.Range("B9").Value = 1600
$ws.Range("C9").Value = 1550

$ws.Range("B10").Value = 1470
$ws.Range("C10").Value = 1475

$ws.Range("B11").Value = 1375
$ws.Range("C11").Value = 1400

$ws.Range("B22").Value = 40491
$ws.Range("C22").Value = 41768

$ws.Range("B23").Value = 26112
$ws.Range("B24").Value = 7750
$ws.Range("B25").Value = 4324
$ws.Range("B26").Value = 3829
$ws.Range("B27").Value = 2371
$ws.Range("B28").Value = 1854
$ws.Range("B29").Value = 1911
$ws.Range("B30").Value = 1376
$ws.Range("B31").Value = 1275
$ws.Range("B32").Value = 1252
$ws.Range("B33").Value = 1013

# Update the selection to match the target (active cell moved to G40)
$ws.Range("G40").Select()
